$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"0.6666666666666666"
$ws.Range("G2").Value = [double]"0.04155"
$ws.Range("H2").Value = [double]"0.12465"
$ws.Range("I2").Value = [double]"0.0001466168179836329"
$ws.Range("J2").Value = [double]"0.0001466168179836329"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.3333333333333333"
$ws.Range("M2").Value = [double]"0.002689333333333334"
$ws.Range("N2").Value = [double]"0.008068000000000001"
$ws.Range("O2").Value = [double]"0.03638167388167388"
$ws.Range("P2").Value = [double]"0.03638167388167389"
$ws.Range("Q2").Value = [double]"0.0001117418"
$ws.Range("R2").Value = [double]"0.0010056762"
$ws.Range("S2").Value = [double]"5.334165257449269E-06"
$ws.Range("T2").Value = [double]"5.334165257449271E-06"
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"0.6666666666666666"
$ws.Range("G3").Value = [double]"0.04155"
$ws.Range("H3").Value = [double]"0.12465"
$ws.Range("I3").Value = [double]"0.0001466168179836329"
$ws.Range("J3").Value = [double]"0.0001466168179836329"
$ws.Range("O3").Value = [double]"0.9259334415584415"
$ws.Range("P3").Value = [double]"0.9259334415584416"
$ws.Range("Q3").Value = [double]"0.002843889749999999"
$ws.Range("R3").Value = [double]"0.02559500775"
$ws.Range("S3").Value = [double]"0.0001357574148659328"
$ws.Range("T3").Value = [double]"0.0001357574148659328"
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"0.6666666666666666"
$ws.Range("G4").Value = [double]"0.04155"
$ws.Range("H4").Value = [double]"0.12465"
$ws.Range("I4").Value = [double]"0.0001466168179836329"
$ws.Range("J4").Value = [double]"0.0001466168179836329"
$ws.Range("O4").Value = [double]"0.03768488455988456"
$ws.Range("P4").Value = [double]"0.03768488455988456"
$ws.Range("Q4").Value = [double]"0.00011574445"
$ws.Range("R4").Value = [double]"0.00104170005"
$ws.Range("S4").Value = [double]"5.525237860250811E-06"
$ws.Range("T4").Value = [double]"5.525237860250811E-06"
$ws.Range("I5").Value = [double]"0.9992428949822291"
$ws.Range("J5").Value = [double]"0.9992428949822291"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.002689333333333334"
$ws.Range("N5").Value = [double]"0.008068000000000001"
$ws.Range("O5").Value = [double]"0.03638167388167388"
$ws.Range("P5").Value = [double]"0.03638167388167389"
$ws.Range("Q5").Value = [double]"0.7615579253328889"
$ws.Range("R5").Value = [double]"6.854021327996001"
$ws.Range("S5").Value = [double]"0.03635412913382316"
$ws.Range("T5").Value = [double]"0.03635412913382317"
$ws.Range("I6").Value = [double]"0.9992428949822291"
$ws.Range("J6").Value = [double]"0.9992428949822291"
$ws.Range("O6").Value = [double]"0.9259334415584415"
$ws.Range("P6").Value = [double]"0.9259334415584416"
$ws.Range("S6").Value = [double]"0.9252324127037157"
$ws.Range("T6").Value = [double]"0.9252324127037158"
$ws.Range("I7").Value = [double]"0.9992428949822291"
$ws.Range("J7").Value = [double]"0.9992428949822291"
$ws.Range("O7").Value = [double]"0.03768488455988456"
$ws.Range("P7").Value = [double]"0.03768488455988456"
$ws.Range("S7").Value = [double]"0.03765635314469016"
$ws.Range("T7").Value = [double]"0.03765635314469016"
$ws.Range("I8").Value = [double]"0.0006104881997874136"
$ws.Range("J8").Value = [double]"0.0006104881997874135"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.002689333333333334"
$ws.Range("N8").Value = [double]"0.008068000000000001"
$ws.Range("O8").Value = [double]"0.03638167388167388"
$ws.Range("P8").Value = [double]"0.03638167388167389"
$ws.Range("Q8").Value = [double]"0.0004652743884444444"
$ws.Range("R8").Value = [double]"0.004187469496"
$ws.Range("S8").Value = [double]"2.221058259327585E-05"
$ws.Range("T8").Value = [double]"2.221058259327585E-05"
$ws.Range("I9").Value = [double]"0.0006104881997874136"
$ws.Range("J9").Value = [double]"0.0006104881997874135"
$ws.Range("O9").Value = [double]"0.9259334415584415"
$ws.Range("P9").Value = [double]"0.9259334415584416"
$ws.Range("S9").Value = [double]"0.0005652714398599772"
$ws.Range("T9").Value = [double]"0.0005652714398599772"
$ws.Range("I10").Value = [double]"0.0006104881997874136"
$ws.Range("J10").Value = [double]"0.0006104881997874135"
$ws.Range("O10").Value = [double]"0.03768488455988456"
$ws.Range("P10").Value = [double]"0.03768488455988456"
$ws.Range("S10").Value = [double]"2.300617733416042E-05"
$ws.Range("T10").Value = [double]"2.300617733416042E-05"
